# Apply "update schema to test" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in the three newly-added shared-string-backed cells (schema
#    columns UserId / SessionId / OTP) on rows 8 and 9.
$ws.Range("A8").Value = "UserId"
$ws.Range("C8").Value = "SessionId"
$ws.Range("A9").Value = "OTP"

# 2. Change the highlight fill color used by the data rows (A2:C9) from
#    green to yellow. (RGB(255,255,0) = 255 + 255*256 + 0*65536 = 65535)
$ws.Range("A2:C9").Interior.Color = 65535

# 3. Move the active selection to A6.
$ws.Range("A6").Select()
